$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The inflow.csv block (rows 10-17) loses the "OGM_don" (dissolved organic
# nitrogen) row that used to sit at row 14, and gains a new "OGM_doc"
# (dissolved organic carbon) row at the end of the block (row 17).
# Concretely: rows 15-17 (NIT_nit, NIT_amm, PHS_frp) each shift up by one
# row, and a brand new OGM_doc row is written at row 17 using the value
# that used to belong to the old row 14 (92.6015) together with a new
# description.

# Row 14: was OGM_don -> becomes NIT_nit (copy old row 15 content up)
$ws.Range("B14").Value = "NIT_nit"
$ws.Range("C14").Value = "mmol N/m3"
$ws.Range("E14").Value = 926.01530000000002
$ws.Range("F14").Value = "nitrate concentration in inflow water"

# Row 15: was NIT_nit -> becomes NIT_amm (copy old row 16 content up)
$ws.Range("B15").Value = "NIT_amm"
$ws.Range("C15").Value = "mmol N/m3"
$ws.Range("E15").Value = 14.084899999999999
$ws.Range("F15").Value = "ammonium concentration in inflow water"

# Row 16: was NIT_amm -> becomes PHS_frp (copy old row 17 content up)
$ws.Range("B16").Value = "PHS_frp"
$ws.Range("C16").Value = "mmol P/m3"
$ws.Range("E16").Value = 0.79179641999999995
$ws.Range("F16").Value = "filterable reactive phosphorus concentration in inflow water"
$ws.Range("G16").Value = "convert from mg/L by multiplying by 32.29"

# Row 17: was PHS_frp -> becomes the new OGM_doc row
$ws.Range("B17").Value = "OGM_doc"
$ws.Range("C17").Value = "mmol N/m3"
$ws.Range("E17").Value = 92.601500000000001
$ws.Range("F17").Value = "dissolved organic carbon concentration in inflow water"
$ws.Range("G17").ClearContents()
$ws.Range("H17").ClearContents()

$wb.Save()
